$d = $word.ActiveDocument

# The document currently ends with a single paragraph containing "@Transactional"
# followed by the hidden "_GoBack" bookmark. We replace that whole paragraph
# (including its paragraph mark) with the expanded content: the same
# "@Transactional" paragraph (now a plain paragraph, no longer carrying the
# paragraph-mark run formatting since it is no longer the last paragraph),
# a blank paragraph, a new "二级回复" sub-heading, a short note, and a final
# paragraph documenting a breakpoint, with the "_GoBack" bookmark following
# the most recently typed text.

$last = $d.Paragraphs.Last
$target = $last.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>@T</w:t></w:r><w:r><w:t>ransactional</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="2"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:lastRenderedPageBreak/><w:t>二级回复</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>弄个破js弄了半天。。。。。</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>断点：comment表缺comment</w:t></w:r><w:r><w:t>Count</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>属性</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>
'@

$target.InsertXML($xml)
